# "Generate Report for Handback"
#
# The handback-status report is regenerated: the two tracked files
# (f6521733-...md and cf2b8dec-...md) swap table positions (cf2b8dec now
# sorts first / f6521733 second) and cf2b8dec's row picks up the data for
# its own handback that just completed (new handoff/handback timestamps,
# new target xlf names, status flips from "Ready for handoff" to
# "Handed back: in sync with en-US", and its stale "version mismatch"
# error clears).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: File Name / Path And Name / ... / zh-cn / de-de / Latest HO Xliff Generate Date
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

# Row 2 now describes cf2b8dec (used to describe f6521733)
$overview.Range("A2").Value = "cf2b8dec-1757-4cdd-8203-e364c866ae09.md"
$overview.Range("B2").Value = "e2e\cf2b8dec-1757-4cdd-8203-e364c866ae09.md"
$overview.Range("G2").Value = "2016-10-27 09:33:38"

# Row 3 now describes f6521733 (used to describe cf2b8dec)
$overview.Range("A3").Value = "f6521733-974e-4003-a153-c5fee60ff6b8.md"
$overview.Range("B3").Value = "e2e\f6521733-974e-4003-a153-c5fee60ff6b8.md"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"
$overview.Range("G3").Value = "2016-10-27 09:29:54"

# ---------------------------------------------------------------------
# Per-locale detail sheets: zh-cn and de-de share the same column layout
#   A Source File Name, C Status, G Latest Handoff File, H Latest Handoff Datetime,
#   I Latest Target File, J Latest Handback File, K Latest Handback DateTime,
#   P Error Detail
# ---------------------------------------------------------------------
foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $suffix = if ($sheetName -eq "zh-cn") { "zh-cn" } else { "de-de" }

    # Row 2 now describes cf2b8dec's freshly-completed handback
    $ws.Range("A2").Value = "cf2b8dec-1757-4cdd-8203-e364c866ae09.md"
    $ws.Range("I2").Value = "cf2b8dec-1757-4cdd-8203-e364c866ae09.md"
    $ws.Range("G2").Value = "cf2b8dec-1757-4cdd-8203-e364c866ae09.5e6bc4810939fd2b6d46b19feab70b58cca069c3.$suffix.xlf"
    $ws.Range("J2").Value = "cf2b8dec-1757-4cdd-8203-e364c866ae09.5e6bc4810939fd2b6d46b19feab70b58cca069c3.$suffix.xlf"

    # Row 3 now describes f6521733
    $ws.Range("A3").Value = "f6521733-974e-4003-a153-c5fee60ff6b8.md"
    $ws.Range("I3").Value = "f6521733-974e-4003-a153-c5fee60ff6b8.md"
    $ws.Range("G3").Value = "f6521733-974e-4003-a153-c5fee60ff6b8.7ab767b4503268295ffbe4ce1ca5a2bd74d23c7e.$suffix.xlf"
    $ws.Range("J3").Value = "f6521733-974e-4003-a153-c5fee60ff6b8.7ab767b4503268295ffbe4ce1ca5a2bd74d23c7e.$suffix.xlf"

    # Row 3's status catches up to "handed back" and its stale-version error clears
    $ws.Range("C3").Value = "Handed back: in sync with en-US"
    $ws.Range("P3").Value = ""
}

# zh-cn specific timestamps
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H2").Value = "2016-10-27 09:33:24"
$zhcn.Range("K2").Value = "2016-10-27 09:34:16"
$zhcn.Range("H3").Value = "2016-10-27 09:29:38"

# de-de specific timestamps
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H2").Value = "2016-10-27 09:33:38"
$dede.Range("K2").Value = "2016-10-27 09:34:33"
$dede.Range("H3").Value = "2016-10-27 09:29:54"
